$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row r, Coin (B), Link (C), Price (D), Volume(1h) (E)
$rows = @(
    @(2, "Bitcoin", "https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc", "71.233.00", "  +2.72%  "),
    @(3, "Ethereum", "https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth", "3.585.05", "  +1.84%  "),
    @(4, "TetherUSD", "https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt", "1.00", "  +0.02%  "),
    @(5, "BNB", "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb", "584.03", "  +2.44%  "),
    @(6, "Solana", "https://coinranking.com/coin/zNZHO_Sjf+solana-sol", "186.57", "  +2.17%  "),
    @(7, "LidoStakedEther", "https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth", "3.573.63", "  +1.70%  "),
    @(8, "XRP", "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp", "0.622", "  +1.32%  "),
    @(9, "USDC", "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc", "1.00", "  +0.06%  "),
    @(10, "Dogecoin", "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge", "0.215", "  +16.30%  "),
    @(11, "Cardano", "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada", "0.652", "  +1.95%  "),
    @(12, "Avalanche", "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax", "54.50", "  +1.77%  "),
    @(13, "ShibaInu", "https://coinranking.com/coin/xz24e0BjL+shibainu-shib", "0.0000318", "  +6.53%  "),
    @(14, "Polkadot", "https://coinranking.com/coin/25W7FG7om+polkadot-dot", "9.53", "  +0.85%  "),
    @(15, "WrappedliquidstakedEther2.0", "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth", "4.022.04", "  -1.77%  "),
    @(16, "WrappedBTC", "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc", "71.240.17", "  +2.78%  "),
    @(17, "Chainlink", "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link", "19.24", "  -0.15%  "),
    @(18, "WrappedEther", "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth", "3.584.44", "  +1.38%  "),
    @(19, "Uniswap", "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni", "12.34", "  -0.11%  "),
    @(20, "BitcoinCash", "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch", "564.79", "  +4.39%  "),
    @(21, "TRON", "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx", "0.120", "  +0.59%  "),
    @(22, "Polygon", "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic", "1.01", "  -1.75%  "),
    @(23, "InternetComputer(DFINITY)", "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp", "17.56", "  -10.04%  "),
    @(24, "Toncoin", "https://coinranking.com/coin/67YlI0K1b+toncoin-ton", "5.10", "  +3.65%  "),
    @(25, "PancakeSwap", "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake", "4.58", "  +5.24%  "),
    @(26, "Litecoin", "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc", "94.84", "  +1.28%  "),
    @(27, "RenderToken", "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr", "11.26", "  +1.00%  "),
    @(28, "ImmutableX", "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx", "2.94", "  +1.14%  "),
    @(29, "Filecoin", "https://coinranking.com/coin/ymQub4fuB+filecoin-fil", "9.14", "  +0.48%  "),
    @(30, "EthereumClassic", "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc", "32.54", "  +2.76%  "),
    @(31, "NEARProtocol", "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near", "7.25", "  -1.62%  "),
    @(32, "Cosmos", "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom", "12.29", "  -1.77%  "),
    @(33, "Hedera", "https://coinranking.com/coin/jad286TjB+hedera-hbar", "0.115", "  +1.03%  "),
    @(34, "OKB", "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb", "64.22", "  -1.08%  "),
    @(35, "Fetch.AI", "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet", "3.34", "  +6.53%  "),
    @(36, "Bittensor", "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao", "547.65", "  -4.42%  "),
    @(37, "TheGraph", "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt", "0.416", "  +4.41%  "),
    @(38, "PEPE", "https://coinranking.com/coin/03WI8NQPF+pepe-pepe", "0.0₃0806", "  +5.62%  "),
    @(39, "InjectiveProtocol", "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj", "37.61", "  -1.10%  "),
    @(40, "Dai", "https://coinranking.com/coin/MoTuySvg7+dai-dai", "1.00", "  +0.16%  "),
    @(41, "dogwifhat", "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif", "3.32", "  +6.60%  "),
    @(42, "Maker", "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr", "3.508.48", "  +11.62%  "),
    @(43, "Stacks", "https://coinranking.com/coin/mMPrMcB7+stacks-stx", "3.47", "  +2.92%  "),
    @(44, "Kaspa", "https://coinranking.com/coin/V8GxkwWow+kaspa-kas", "0.136", "  +1.48%  "),
    @(45, "VeChain", "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet", "0.0446", "  +0.63%  "),
    @(46, "ApeXProtocol", "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex", "3.48", "  -1.25%  "),
    @(47, "ThetaToken", "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta", "2.94", "  -0.44%  "),
    @(48, "THORChain", "https://coinranking.com/coin/ybmU-kKU+thorchain-rune", "9.40", "  +2.15%  "),
    @(49, "Stellar", "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm", "0.137", "  +2.42%  "),
    @(50, "FirstDigitalUSD", "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd", "0.999", "  -0.35%  "),
    @(51, "OceanProtocol", "https://coinranking.com/coin/aAKLSV5-0+oceanprotocol-ocean", "1.45", "  +2.14%  ")
)

foreach ($row in $rows) {
    $r = $row[0]
    $coin = $row[1]
    $link = $row[2]
    $price = $row[3]
    $volume = $row[4]

    $ws.Range("B$r").Value = $coin
    $ws.Range("C$r").Value = $link

    $ws.Range("D$r").NumberFormat = "@"
    $ws.Range("D$r").Value = $price

    $ws.Range("E$r").NumberFormat = "@"
    $ws.Range("E$r").Value = $volume
}
